$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.113.81"
$ws.Range("E2").Value = "  +1.62%  "

$ws.Range("D3").Value = "3.213.95"
$ws.Range("E3").Value = "  +1.06%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.211.94"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.510"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").Value = "3.741.55"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.233.14"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.86%  "

$ws.Range("D18").Value = "3.224.21"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.23"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.01%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "

$ws.Range("E30").Value = "  +10.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0916"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "485.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0423"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.297"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.05%  "

$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.31%  "

$ws.Range("D45").Value = "2.963.11"
$ws.Range("E45").Value = "  -3.46%  "

$ws.Range("D46").Value = "0.0₃0646"
$ws.Range("E46").Value = "  +4.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.116"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.46%  "
